$p = $ppt.ActivePresentation

# Remove the "Let's See it in action! / Student Demonstration" slide
# (SlideID 319) - it was deleted from the deck. Locate it by SlideID so
# the script is robust to index shifts.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 319) {
        $s.Delete()
    }
}

# The deck was re-saved a few days later - the cached "today" date shown
# on the handout master and notes master footers moved from 3/30/21 to
# 4/2/21.
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "4/2/21"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "4/2/21"
